# Hazard extraction accuracy edit
# - Rewrites the "Hazard-focused" summary table (rows 2-18) with the
#   updated hazard/action/negation word lists produced after improving
#   the hazard extraction accuracy.
# - Updates the saved selection / active sheet to match what was left
#   selected when the workbook was saved.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hazard-focused")
$ws2 = $wb.Worksheets.Item("topic-focused")

$ws1.Cells.Item(2,1).Value = "highway, traffic, road, travel, interstate"
$ws1.Cells.Item(2,2).Value = "closure, remain, remains, close, block, impact, access, limit, limited"
$ws1.Cells.Item(2,3).Value = "reopen, open, lift"
$ws1.Cells.Item(2,4).Value = 6
$ws1.Cells.Item(2,10).Value = "Mission"
$ws1.Cells.Item(2,11).Value = "Traffic"

$ws1.Cells.Item(3,1).Value = "jurisdiction, team, command, organization, type"
$ws1.Cells.Item(3,2).Value = "involve, transition, transfer"
$ws1.Cells.Item(3,3).Value = $null
$ws1.Cells.Item(3,4).Value = "4, 72"
$ws1.Cells.Item(3,10).Value = "Mission"
$ws1.Cells.Item(3,11).Value = "Command Transitions"

$ws1.Cells.Item(4,1).Value = "evacuation, evacuate, threaten, threat"
$ws1.Cells.Item(4,2).Value = "resident, residence, level, notice, community, structure, subdivision, mandatory, order, effect, remain, continue, issue"
$ws1.Cells.Item(4,3).Value = "lift, return, reopen, open, reduce"
$ws1.Cells.Item(4,4).Value = "25, 45, 146"
$ws1.Cells.Item(4,10).Value = "Mission"
$ws1.Cells.Item(4,11).Value = "Evacuations"

$ws1.Cells.Item(5,1).Value = "mapping, map"
$ws1.Cells.Item(5,2).Value = "reflect, accurate, adjustment, change, reflect, inaccurate"
$ws1.Cells.Item(5,3).Value = $null
$ws1.Cells.Item(5,4).Value = 135
$ws1.Cells.Item(5,10).Value = "Mission"
$ws1.Cells.Item(5,11).Value = "Inaccurate Mapping"

$ws1.Cells.Item(6,1).Value = "aircraft, heli, helicopter, aerial, tanker, copter, ground"
$ws1.Cells.Item(6,2).Value = "suspend, smoke, hazard, windy, wind, suspendsion, mechanical, problem, due"
$ws1.Cells.Item(6,3).Value = "resume, drop, lack, lift"
$ws1.Cells.Item(6,4).Value = $null
$ws1.Cells.Item(6,10).Value = "Mission"
$ws1.Cells.Item(6,11).Value = "Aerial Grounding"

$ws1.Cells.Item(7,1).Value = "unstaffed, resource, support, crew, aircraft, helicopter, engine, staff"
$ws1.Cells.Item(7,2).Value = "lack, need, shortage, minimal, share, necessary, limited, limit, fatigue"
$ws1.Cells.Item(7,3).Value = "release, demob, demobilization, demobilize, progress"
$ws1.Cells.Item(7,4).Value = "99, 152"
$ws1.Cells.Item(7,10).Value = "Mission"
$ws1.Cells.Item(7,11).Value = "Resource Issues"

$ws1.Cells.Item(8,1).Value = "injury, hospital, injured, accident, treatment, laceration, firefighter, treat"
$ws1.Cells.Item(8,2).Value = "minor, report, transport, heat, shoulder, ankle, medical, release"
$ws1.Cells.Item(8,3).Value = $null
$ws1.Cells.Item(8,4).Value = 34
$ws1.Cells.Item(8,10).Value = "Mission"
$ws1.Cells.Item(8,11).Value = "Injuries"

$ws1.Cells.Item(9,1).Value = "impact, concern, site, nation"
$ws1.Cells.Item(9,2).Value = "political, social, adjacent, community, cultural, tribal, monument, archaeological, heritage"
$ws1.Cells.Item(9,3).Value = "smoke"
$ws1.Cells.Item(9,4).Value = 5
$ws1.Cells.Item(9,10).Value = "Wildland Urban Interface"
$ws1.Cells.Item(9,11).Value = "Cultural Resources"

$ws1.Cells.Item(10,1).Value = "cattle, buffalo, allotment, ranch, sheep, livestock"
$ws1.Cells.Item(10,2).Value = "grazing, pasture, threaten, concern, risk, threat, private, area, evacuate, evacuation, order"
$ws1.Cells.Item(10,3).Value = $null
$ws1.Cells.Item(10,4).Value = 23
$ws1.Cells.Item(10,10).Value = "Wildland Urban Interface"
$ws1.Cells.Item(10,11).Value = "Livestock"

$ws1.Cells.Item(11,1).Value = "violation, notification, respond"
$ws1.Cells.Item(11,2).Value = "law, patrol"
$ws1.Cells.Item(11,3).Value = "reverse"
$ws1.Cells.Item(11,4).Value = 21
$ws1.Cells.Item(11,10).Value = "Wildland Urban Interface"
$ws1.Cells.Item(11,11).Value = "Law Violations"

$ws1.Cells.Item(12,1).Value = "military, unexploded"
$ws1.Cells.Item(12,2).Value = "training, present, ordinance, proximity, activity, active, base, area"
$ws1.Cells.Item(12,3).Value = $null
$ws1.Cells.Item(12,4).Value = 52
$ws1.Cells.Item(12,10).Value = "Wildland Urban Interface"
$ws1.Cells.Item(12,11).Value = "Military Base"

$ws1.Cells.Item(13,1).Value = "infrastructure, utility, powerline, water, electric, pipeline, powerlines, watershed, pole, power, gas"
$ws1.Cells.Item(13,2).Value = "concern, near, hazard, critical, threaten, threat, off"
$ws1.Cells.Item(13,3).Value = "restore, tender, diminished"
$ws1.Cells.Item(13,4).Value = "70, 71"
$ws1.Cells.Item(13,10).Value = "Wildland Urban Interface"
$ws1.Cells.Item(13,11).Value = "Infrastructure"

$ws1.Cells.Item(14,1).Value = "weather, behavior, wind, thunderstorm, storm, gusty, lightning, flag"
$ws1.Cells.Item(14,2).Value = "unpredictable, extreme, erratic, strong, red, warning, warn"
$ws1.Cells.Item(14,3).Value = $null
$ws1.Cells.Item(14,4).Value = 76
$ws1.Cells.Item(14,10).Value = "Environmental"
$ws1.Cells.Item(14,11).Value = "Extreme Weather"

$ws1.Cells.Item(15,1).Value = "specie, habitat, animal, plant, conservation"
$ws1.Cells.Item(15,2).Value = "threaten, endanger, threat, sensitive, risk, loss, impact"
$ws1.Cells.Item(15,3).Value = $null
$ws1.Cells.Item(15,4).Value = 80
$ws1.Cells.Item(15,10).Value = "Environmental"
$ws1.Cells.Item(15,11).Value = "Ecological"

$ws1.Cells.Item(16,1).Value = "terrain, rollout, snag, steep, debris, access"
$ws1.Cells.Item(16,2).Value = "concern, hazardous, pose, heavy, rugged, difficult, steep, narrow"
$ws1.Cells.Item(16,3).Value = $null
$ws1.Cells.Item(16,4).Value = "20, 33, 89, 120"
$ws1.Cells.Item(16,10).Value = "Environmental"
$ws1.Cells.Item(16,11).Value = "Hazardous Terrain"

$ws1.Cells.Item(17,1).Value = "flood, flashflood"
$ws1.Cells.Item(17,2).Value = "flash, risk, potential, chance"
$ws1.Cells.Item(17,3).Value = $null
$ws1.Cells.Item(17,4).Value = 128
$ws1.Cells.Item(17,10).Value = "Environmental"
$ws1.Cells.Item(17,11).Value = "Floods"

$ws1.Cells.Item(18,1).Value = "humidity, moisture, hot, drought"
$ws1.Cells.Item(18,2).Value = "low, dry, prolong"
$ws1.Cells.Item(18,3).Value = $null
$ws1.Cells.Item(18,4).Value = 144
$ws1.Cells.Item(18,10).Value = "Environmental"
$ws1.Cells.Item(18,11).Value = "Dry Weather"

# --- view / selection state -------------------------------------------------
# Leave a selection behind on the topic-focused sheet (as last edited),
# then switch back to Hazard-focused as the active/visible tab.
$ws2.Select()
$ws2.Range("I141").Select()

$ws1.Select()
$ws1.Range("K24").Select()
